$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row (row 11) - correct answer marking value
$ws.Range("B11").Value = 5

# Update the "Total" row (row 12) - total marks obtained and the "correct/total" display
$ws.Range("B12").Value = 135
$ws.Range("E12").Value = "135/140"
